# Data correction pass ("Updated by the robot"):
# Updates canton VD (column X) and national total CH (column AB) figures
# on the Cases sheet, plus matching corrections on the Hospitalized and
# ICU sheets (columns X, L and AB).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("X35").Value = 182
$ws.Range("AB35").Value = 1276
$ws.Range("X36").Value = 185
$ws.Range("AB36").Value = 1318
$ws.Range("X37").Value = 191
$ws.Range("AB37").Value = 1370
$ws.Range("X38").Value = 204
$ws.Range("AB38").Value = 1455
$ws.Range("X39").Value = 220
$ws.Range("AB39").Value = 1546
$ws.Range("X40").Value = 235
$ws.Range("AB40").Value = 1649
$ws.Range("X41").Value = 247
$ws.Range("AB41").Value = 1742
$ws.Range("X42").Value = 260
$ws.Range("AB42").Value = 1839
$ws.Range("X43").Value = 268
$ws.Range("AB43").Value = 1893
$ws.Range("X44").Value = 273
$ws.Range("AB44").Value = 1953
$ws.Range("X45").Value = 286
$ws.Range("AB45").Value = 2062
$ws.Range("X46").Value = 305
$ws.Range("AB46").Value = 2175
$ws.Range("X47").Value = 321
$ws.Range("AB47").Value = 2315
$ws.Range("X48").Value = 334
$ws.Range("AB48").Value = 2420
$ws.Range("X49").Value = 350
$ws.Range("AB49").Value = 2534
$ws.Range("X50").Value = 361
$ws.Range("AB50").Value = 2594
$ws.Range("X51").Value = 367
$ws.Range("AB51").Value = 2642
$ws.Range("X52").Value = 380
$ws.Range("AB52").Value = 2766
$ws.Range("X53").Value = 386
$ws.Range("AB53").Value = 2909
$ws.Range("X54").Value = 397
$ws.Range("AB54").Value = 3037
$ws.Range("X59").Value = 462
$ws.Range("AB59").Value = 3630
$ws.Range("X60").Value = 489
$ws.Range("AB60").Value = 3785
$ws.Range("X61").Value = 505
$ws.Range("AB61").Value = 3935
$ws.Range("AB62").Value = 4044
$ws.Range("L62").Value = 55

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("X6").Value = 2
$ws.Range("AB6").Value = -26
$ws.Range("X7").Value = 3
$ws.Range("AB7").Value = -27
$ws.Range("X8").Value = -1
$ws.Range("AB8").Value = -37
$ws.Range("X9").Value = -1
$ws.Range("AB9").Value = -40
$ws.Range("X10").Value = 0
$ws.Range("AB10").Value = -54
$ws.Range("X11").Value = -1
$ws.Range("AB11").Value = -53
$ws.Range("X12").Value = -2
$ws.Range("AB12").Value = -56
$ws.Range("X13").Value = -2
$ws.Range("AB13").Value = -65
$ws.Range("X14").Value = 0
$ws.Range("AB14").Value = -58
$ws.Range("X15").Value = -1
$ws.Range("AB15").Value = -63
$ws.Range("X16").Value = -1
$ws.Range("AB16").Value = -63
$ws.Range("X17").Value = -2
$ws.Range("AB17").Value = -75
$ws.Range("X18").Value = -5
$ws.Range("AB18").Value = -80
$ws.Range("X19").Value = -7
$ws.Range("AB19").Value = -77
$ws.Range("X20").Value = -8
$ws.Range("AB20").Value = -79
$ws.Range("X21").Value = -6
$ws.Range("AB21").Value = -81
$ws.Range("X22").Value = -6
$ws.Range("AB22").Value = -82
$ws.Range("X23").Value = -5
$ws.Range("AB23").Value = -80
$ws.Range("X24").Value = -9
$ws.Range("AB24").Value = -78
$ws.Range("X25").Value = -11
$ws.Range("AB25").Value = -83
$ws.Range("X26").Value = -12
$ws.Range("AB26").Value = -83
$ws.Range("X27").Value = -13
$ws.Range("AB27").Value = -87
$ws.Range("X28").Value = -12
$ws.Range("AB28").Value = -85
$ws.Range("X29").Value = -12
$ws.Range("AB29").Value = -86
$ws.Range("X30").Value = -12
$ws.Range("AB30").Value = -87
$ws.Range("X31").Value = -11
$ws.Range("AB31").Value = -79
$ws.Range("X32").Value = -12
$ws.Range("AB32").Value = -77
$ws.Range("X33").Value = -13
$ws.Range("AB33").Value = -77
$ws.Range("X34").Value = -12
$ws.Range("AB34").Value = -84
$ws.Range("X35").Value = -13
$ws.Range("AB35").Value = -85
$ws.Range("X36").Value = -14
$ws.Range("AB36").Value = -85
$ws.Range("X37").Value = -12
$ws.Range("AB37").Value = -81
$ws.Range("X38").Value = -18
$ws.Range("AB38").Value = -83
$ws.Range("X39").Value = -18
$ws.Range("AB39").Value = -83
$ws.Range("X40").Value = -18
$ws.Range("AB40").Value = -90
$ws.Range("X41").Value = -18
$ws.Range("AB41").Value = -88
$ws.Range("X42").Value = -16
$ws.Range("AB42").Value = -83
$ws.Range("X43").Value = -15
$ws.Range("AB43").Value = -85
$ws.Range("X44").Value = -15
$ws.Range("AB44").Value = -77
$ws.Range("X45").Value = -18
$ws.Range("AB45").Value = -64
$ws.Range("X46").Value = -20
$ws.Range("AB46").Value = -69
$ws.Range("X47").Value = -20
$ws.Range("AB47").Value = -72
$ws.Range("X48").Value = -18
$ws.Range("AB48").Value = -73
$ws.Range("X49").Value = -16
$ws.Range("AB49").Value = -61
$ws.Range("X50").Value = -13
$ws.Range("AB50").Value = -62
$ws.Range("X51").Value = -11
$ws.Range("AB51").Value = -58
$ws.Range("X52").Value = -11
$ws.Range("AB52").Value = -58
$ws.Range("X53").Value = -12
$ws.Range("AB53").Value = -58
$ws.Range("X54").Value = -14
$ws.Range("AB54").Value = -67
$ws.Range("X55").Value = -13
$ws.Range("AB55").Value = -66
$ws.Range("X59").Value = -11
$ws.Range("AB59").Value = -61
$ws.Range("X60").Value = -14
$ws.Range("AB60").Value = -66
$ws.Range("X61").Value = -11
$ws.Range("AB61").Value = -55
$ws.Range("AB62").Value = -49
$ws.Range("L62").Value = -3

$ws = $wb.Worksheets.Item("ICU")
$ws.Range("X60").Value = -4
$ws.Range("AB60").Value = -9
$ws.Range("X61").Value = -3
$ws.Range("AB61").Value = -9
$ws.Range("AB62").Value = -6
$ws.Range("L62").Value = -1

